$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix accented question/field text (Como -> Cómo, Que -> Qué, Donde -> Dónde, etc.) ---
$ws.Range("B2").Value = "¿Qué número de la pokedex es #pokemon#?, ¿Qué entrada es #pokemon# en la pokedex?, ¿Cómo se identifica a #pokemon#?, ¿Cuál es el identificador a #pokemon#?, ¿En qué página esta #pokemon#?"

$ws.Range("A3").Value = "Descripcion"
$ws.Range("B3").Value = "¿Cuál es la descripción de #pokemon#?, ¿Cómo es #pokemon#?, ¿Cómo describirías a #pokemon#?"

$ws.Range("B5").Value = "¿A que es débil #pokemon#?, ¿Qué le hace mucho daño a #pokemon#?, ¿Qué tiene que evitar #pokemon#?, ¿Qué debilita o hiere a #pokemon#?, ¿Qué le afecta a #pokemon#?, ¿Cómo debilitar a pikachu?"

$ws.Range("B6").Value = "¿Cuánto pesa #pokemon#?, ¿Cómo de pesado es #pokemon#?"

$ws.Range("B7").Value = "¿Cuánto mide #pokemon#?, ¿Cómo de alto es #pokemon#?, ¿Cómo de grande es #pokemon#?"

$ws.Range("A11").Value = "Generacion"

$ws.Range("A14").Value = "Obtencion"
$ws.Range("B14").Value = "¿Cómo se obtiene a #pokemon#?, ¿Cómo se consigue a #pokemon#?, ¿Dónde se consigue a #pokemon#?, ¿Cómo se captura a #pokemon#?, ¿Dónde se captura a #pokemon#?"

$ws.Range("B15").Value = "¿Cuál es el ratio de captura de #pokemon#?, ¿Es fácil capturar a #pokemon#?, ¿Es complicado capturar a #pokemon#?, ¿Cómo de probable es capturar a #pokemon#?"

$ws.Range("B16").Value = "¿Qué movimientos aprende #pokemon#?, ¿Qué aprende #pokemon#?, ¿Qué ataques puede hacer #pokemon#?, ¿Con que puede atacar #pokemon#?"

# --- Clear the now-unused bold / bold+underline emphasis from the section-label cells ---
$ws.Range("A2").ClearFormats()
$ws.Range("A10").ClearFormats()
$ws.Range("A13").ClearFormats()
$ws.Range("A15").ClearFormats()

# --- Move the active selection to match where the author last clicked ---
$ws.Range("B16").Select()
